$d = $word.ActiveDocument

$replacements = @(
    @("892÷2=", "673÷2="),
    @("825÷9=", "862÷8="),
    @("798÷7=", "205÷3="),
    @("994÷6=", "138÷3="),
    @("966÷8=", "735÷4="),
    @("284÷8=", "269÷4="),
    @("382÷2=", "786÷5="),
    @("823÷7=", "188÷2="),
    @("291÷2=", "120÷9="),
    @("548÷6=", "199÷7="),
    @("484÷5=", "201÷2="),
    @("101÷8=", "356÷9="),
    @("858÷2=", "912÷7="),
    @("374÷6=", "856÷4="),
    @("567÷9=", "827÷3="),
    @("100÷5=", "259÷3="),
    @("866÷9=", "188÷4="),
    @("520÷8=", "140÷9="),
    @("415÷8=", "467÷8="),
    @("752÷5=", "845÷4="),
    @("362÷9=", "430÷4="),
    @("449÷3=", "134÷7="),
    @("505÷8=", "197÷6="),
    @("442÷5=", "171÷3="),
    @("388÷9=", "327÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
